$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 97: fully entered first (Name, Website, Type, Loc, Info) ---
$ws.Range("A97").Value = "OASAS Naloxone Training Calendar"
$ws.Range("D97").Value = "https://oasas.ny.gov/keywords/naloxone"
$ws.Range("G97").Value = "AOD"
$ws.Range("H97").Value = "Off"
$ws.Range("I97").Value = "Register for virtual or in person Naloxone (NARCAN) trainings. If you attend virtually, a free Naloxone kit will be mailed to your address."

# --- Names for rows 98-100 entered next (down column A) ---
$ws.Range("A98").Value = "SMART Recovery Meetings"
$ws.Range("A99").Value = "Refuge Recovery Meetings"
$ws.Range("A100").Value = "AA Meetings"

# --- Row 98 remaining fields ---
$ws.Range("D98").Value = "https://meetings.smartrecovery.org/meetings/"
$ws.Range("G98").Value = "AOD"
$ws.Range("H98").Value = "Off"
$ws.Range("I98").Value = "Find a SMART Reocvery Meeting. SMART stands for Self Management and Recovery Training. SMART Recovery helps people recover from addictive and problematic behaviors, using a self-empowering and evidence-informed program."

# --- Row 99 remaining fields ---
$ws.Range("D99").Value = "https://refugerecoverymeetings.org/meetings?tsml-day=any&tsml-type=ONL"
$ws.Range("G99").Value = "AOD"
$ws.Range("H99").Value = "Off"
$ws.Range("I99").Value = "Find a Refuge Recovery Meeting. Refuge Recovery is a systematic method based on Buddhist principles, which integrates scientific, non-theistic, and psychological insight. Viewing addiction as craving in the mind and body, Refuge Recovery shows how a path of meditative awareness can alleviate those desires and ease suffering."

# --- Row 100 remaining fields ---
$ws.Range("D100").Value = "https://www.aa.org/find-aa"
$ws.Range("G100").Value = "AOD"
$ws.Range("H100").Value = "Off"
$ws.Range("I100").Value = "Find an Alcoholic's Anonymous Meeting. Alcoholics Anonymous is a fellowship of men and women who share their experience, strength and hope with each other that they may solve their common problem and help others to recover from alcoholism. The only requirement for membership is a desire to stop drinking."

[void]$ws.Range("H95").Select()
$ws.Application.ActiveWindow.ScrollRow = 68
$ws.Application.ActiveWindow.ScrollColumn = 7
